$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "312.60"
Set-TextValue $ws.Cells.Item(2, 5) "3.16%"
Set-TextValue $ws.Cells.Item(2, 7) "6"
Set-TextValue $ws.Cells.Item(3, 4) "38.81"
Set-TextValue $ws.Cells.Item(3, 5) "8.69%"
Set-TextValue $ws.Cells.Item(3, 7) "6"
Set-TextValue $ws.Cells.Item(4, 4) "5.112"
Set-TextValue $ws.Cells.Item(4, 5) "1.75%"
Set-TextValue $ws.Cells.Item(4, 7) "6"
Set-TextValue $ws.Cells.Item(5, 4) "0.08187"
Set-TextValue $ws.Cells.Item(5, 5) "3.13%"
Set-TextValue $ws.Cells.Item(5, 7) "6"
Set-TextValue $ws.Cells.Item(6, 4) "2.011"
Set-TextValue $ws.Cells.Item(6, 5) "7.78%"
Set-TextValue $ws.Cells.Item(6, 7) "6"
Set-TextValue $ws.Cells.Item(7, 4) "7.915"
Set-TextValue $ws.Cells.Item(7, 5) "2.17%"
Set-TextValue $ws.Cells.Item(7, 7) "6"
Set-TextValue $ws.Cells.Item(8, 4) "0.9326"
Set-TextValue $ws.Cells.Item(8, 5) "1.43%"
Set-TextValue $ws.Cells.Item(8, 7) "6"
Set-TextValue $ws.Cells.Item(9, 4) "0.1406"
Set-TextValue $ws.Cells.Item(9, 5) "4.95%"
Set-TextValue $ws.Cells.Item(9, 7) "6"
Set-TextValue $ws.Cells.Item(10, 4) "0.1947"
Set-TextValue $ws.Cells.Item(10, 5) "3.33%"
Set-TextValue $ws.Cells.Item(10, 7) "6"
Set-TextValue $ws.Cells.Item(11, 4) "0.09288"
Set-TextValue $ws.Cells.Item(11, 5) "2.52%"
Set-TextValue $ws.Cells.Item(11, 7) "6"
Set-TextValue $ws.Cells.Item(12, 4) "0.03436"
Set-TextValue $ws.Cells.Item(12, 5) "-0.04%"
Set-TextValue $ws.Cells.Item(12, 7) "6"
Set-TextValue $ws.Cells.Item(13, 4) "0.09850"
Set-TextValue $ws.Cells.Item(13, 5) "0.58%"
Set-TextValue $ws.Cells.Item(13, 7) "6"
Set-TextValue $ws.Cells.Item(14, 4) "0.001408"
Set-TextValue $ws.Cells.Item(14, 5) "0.63%"
Set-TextValue $ws.Cells.Item(14, 7) "6"
Set-TextValue $ws.Cells.Item(15, 4) "0.005990"
Set-TextValue $ws.Cells.Item(15, 5) "-0.96%"
Set-TextValue $ws.Cells.Item(15, 7) "6"
Set-TextValue $ws.Cells.Item(16, 4) "3.652"
Set-TextValue $ws.Cells.Item(16, 5) "-2.24%"
Set-TextValue $ws.Cells.Item(16, 7) "6"
Set-TextValue $ws.Cells.Item(17, 4) "4.185"
Set-TextValue $ws.Cells.Item(17, 5) "2.02%"
Set-TextValue $ws.Cells.Item(17, 7) "6"
Set-TextValue $ws.Cells.Item(18, 4) "3.450"
Set-TextValue $ws.Cells.Item(18, 5) "1.71%"
Set-TextValue $ws.Cells.Item(18, 7) "6"
Set-TextValue $ws.Cells.Item(19, 4) "0.3451"
Set-TextValue $ws.Cells.Item(19, 5) "0.25%"
Set-TextValue $ws.Cells.Item(19, 7) "6"
Set-TextValue $ws.Cells.Item(20, 4) "0.1313"
Set-TextValue $ws.Cells.Item(20, 5) "-1.42%"
Set-TextValue $ws.Cells.Item(20, 7) "6"
Set-TextValue $ws.Cells.Item(21, 4) "4.805"
Set-TextValue $ws.Cells.Item(21, 5) "-7.24%"
Set-TextValue $ws.Cells.Item(21, 7) "6"
Set-TextValue $ws.Cells.Item(22, 4) "0.2453"
Set-TextValue $ws.Cells.Item(22, 5) "2.69%"
Set-TextValue $ws.Cells.Item(22, 7) "6"
Set-TextValue $ws.Cells.Item(23, 4) "0.04474"
Set-TextValue $ws.Cells.Item(23, 5) "1.20%"
Set-TextValue $ws.Cells.Item(23, 7) "6"
Set-TextValue $ws.Cells.Item(24, 4) "0.001239"
Set-TextValue $ws.Cells.Item(24, 5) "2.23%"
Set-TextValue $ws.Cells.Item(24, 7) "6"
Set-TextValue $ws.Cells.Item(25, 5) "-9.78%"
Set-TextValue $ws.Cells.Item(25, 7) "6"
Set-TextValue $ws.Cells.Item(26, 7) "6"
Set-TextValue $ws.Cells.Item(27, 5) "0.19%"
Set-TextValue $ws.Cells.Item(27, 7) "6"
Set-TextValue $ws.Cells.Item(28, 7) "6"
Set-TextValue $ws.Cells.Item(29, 7) "6"
Set-TextValue $ws.Cells.Item(30, 7) "6"
Set-TextValue $ws.Cells.Item(31, 7) "6"
Set-TextValue $ws.Cells.Item(32, 7) "6"
Set-TextValue $ws.Cells.Item(33, 7) "6"
Set-TextValue $ws.Cells.Item(34, 7) "6"
Set-TextValue $ws.Cells.Item(35, 7) "6"
Set-TextValue $ws.Cells.Item(36, 7) "6"
Set-TextValue $ws.Cells.Item(37, 7) "6"
Set-TextValue $ws.Cells.Item(38, 7) "6"
Set-TextValue $ws.Cells.Item(39, 4) "0.02139"
Set-TextValue $ws.Cells.Item(39, 5) "10.36%"
Set-TextValue $ws.Cells.Item(39, 7) "6"
Set-TextValue $ws.Cells.Item(40, 4) "0.05187"
Set-TextValue $ws.Cells.Item(40, 5) "-1.37%"
Set-TextValue $ws.Cells.Item(40, 7) "6"
Set-TextValue $ws.Cells.Item(41, 4) "0.007453"
Set-TextValue $ws.Cells.Item(41, 5) "-1.58%"
Set-TextValue $ws.Cells.Item(41, 7) "6"
Set-TextValue $ws.Cells.Item(42, 4) "0.009993"
Set-TextValue $ws.Cells.Item(42, 5) "-1.11%"
Set-TextValue $ws.Cells.Item(42, 7) "6"
Set-TextValue $ws.Cells.Item(43, 4) "0.1370"
Set-TextValue $ws.Cells.Item(43, 5) "1.71%"
Set-TextValue $ws.Cells.Item(43, 7) "6"
Set-TextValue $ws.Cells.Item(44, 5) "-0.74%"
Set-TextValue $ws.Cells.Item(44, 7) "6"
Set-TextValue $ws.Cells.Item(45, 4) "0.009788"
Set-TextValue $ws.Cells.Item(45, 5) "-3.46%"
Set-TextValue $ws.Cells.Item(45, 7) "6"
Set-TextValue $ws.Cells.Item(46, 5) "3.28%"
Set-TextValue $ws.Cells.Item(46, 7) "6"
Set-TextValue $ws.Cells.Item(47, 5) "0.41%"
Set-TextValue $ws.Cells.Item(47, 7) "6"
Set-TextValue $ws.Cells.Item(48, 7) "6"
Set-TextValue $ws.Cells.Item(49, 5) "-3.14%"
Set-TextValue $ws.Cells.Item(49, 7) "6"
Set-TextValue $ws.Cells.Item(50, 4) "0.00002102"
Set-TextValue $ws.Cells.Item(50, 5) "0.41%"
Set-TextValue $ws.Cells.Item(50, 7) "6"
Set-TextValue $ws.Cells.Item(51, 4) "0.0002002"
Set-TextValue $ws.Cells.Item(51, 5) "0.41%"
Set-TextValue $ws.Cells.Item(51, 7) "6"
